# Updated TestData for Portugal Market
#
# 1. Selection on "Germany" changes to the whole used range (A1:D14).
# 2. A new "Portugal" worksheet is appended after "Swiss", cloned from the
#    "Czech" sheet (same layout/styles), with its own market name / ticket
#    reference and a couple of formatting tweaks (row height, column widths).
# 3. "Portugal" becomes the active/selected sheet (tabSelected moves off
#    "Swiss" onto "Portugal").

$wb = $excel.ActiveWorkbook

# --- 1. Update the stored selection on "Germany" (select whole used range) ---
$germany = $wb.Worksheets.Item("Germany")
$germany.Activate()
$germany.Range("A1:D14").Select()

# --- 2. Clone "Czech" into a brand-new "Portugal" sheet after "Swiss" ---
$czech = $wb.Worksheets.Item("Czech")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$czech.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Portugal"

$portugal = $wb.Worksheets.Item("Portugal")

# Market name / reference cells
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T2433/T2466"

# Row heights for the wrapped "User Story" / "Description" rows
$portugal.Range("A3:D3").RowHeight = 28.8
$portugal.Range("A4:D4").RowHeight = 28.8

# Column widths (engine snaps to 1/6-character increments, so the inputs
# below are pre-compensated so the stored width lands on the closest
# achievable value to the target widths)
$portugal.Columns(1).ColumnWidth = 42.5
$portugal.Columns(2).ColumnWidth = 25.333333333333336
$portugal.Columns(3).ColumnWidth = 13.666666666666666
$portugal.Columns(4).ColumnWidth = 16.666666666666668

# --- 3. Make "Portugal" the active sheet/tab with cursor on B4 ---
$portugal.Activate()
$portugal.Range("B4").Select()
